# Regenerate the K column (column G) values in the save_data sheet.
# The commit re-derives strikeouts (K) from the source data (using
# actual K counts instead of the old "Strike#" proxy), so here we
# directly write the recalculated K values back into column G for
# each data row (rows 2-45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 3
    6  = 1
    7  = 3
    8  = 3
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 3
    18 = 2
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 3
    24 = 4
    25 = 1
    26 = 1
    27 = 4
    28 = 3
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 2
    35 = 0
    36 = 1
    37 = 3
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    42 = 2
    43 = 4
    44 = 3
    45 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
